$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-27 Thursday" "2024-06-28 Friday"

Replace-Text "176×9=1584" "289×9=2601"
Replace-Text "481×3=1443" "899×8=7192"
Replace-Text "273×6=1638" "425×4=1700"
Replace-Text "532×5=2660" "131×7=917"
Replace-Text "929×2=1858" "231×7=1617"

Replace-Text "783×8=6264" "356×2=712"
Replace-Text "664×4=2656" "362×8=2896"
Replace-Text "418×5=2090" "270×4=1080"
Replace-Text "589×7=4123" "553×4=2212"
Replace-Text "815×8=6520" "585×8=4680"

Replace-Text "116×3=348" "752×4=3008"
Replace-Text "336×8=2688" "244×5=1220"
Replace-Text "446×5=2230" "653×5=3265"
Replace-Text "509×3=1527" "683×9=6147"
Replace-Text "965×4=3860" "451×9=4059"

Replace-Text "962×7=6734" "819×4=3276"
Replace-Text "184×6=1104" "161×5=805"
Replace-Text "164×4=656" "187×5=935"
Replace-Text "145×6=870" "776×3=2328"
Replace-Text "106×3=318" "621×8=4968"

Replace-Text "154×6=924" "519×6=3114"
Replace-Text "199×3=597" "842×5=4210"
Replace-Text "982×9=8838" "568×8=4544"
Replace-Text "266×8=2128" "889×4=3556"
Replace-Text "166×8=1328" "254×9=2286"
